# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the "Estado de Cuenta" worker table (rows 16-25) so that period
# "1802" rows come first (rows 16-20), followed by period "1803" rows
# (rows 21-25), and refreshes the "Salario Basico" (column G) values for
# every worker/period combination. "Valor Mora" (column F) and the
# worker identity columns (C = document number, D = name) keep the same
# value per worker, they are simply reordered along with the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered data for rows 16..25 after the update:
# DocNumber, Name, Period, ValorMora, SalarioBasico
$data = @(
    @("8854409",    "NILSON HERRERA PEREZ",            "1802", 38305,  920000),
    @("1032408375", "HECTOR DANIEL GARCIA ABONDANO",   "1802", 166544, 4000000),
    @("52718112",   "ANGELICA MARIA GULFO BASTIDAS",   "1802", 200000, 5000000),
    @("1047471603", "LEONARDO DANIEL ADARRAGA PINTO",  "1802", 52836,  1269000),
    @("73089307",   "GERMAN ANTONIO GARZON GOMEZ",     "1802", 42401,  1018372),
    @("8854409",    "NILSON HERRERA PEREZ",            "1803", 36800,  920000),
    @("1032408375", "HECTOR DANIEL GARCIA ABONDANO",   "1803", 166544, 4000000),
    @("52718112",   "ANGELICA MARIA GULFO BASTIDAS",   "1803", 200000, 5000000),
    @("1047471603", "LEONARDO DANIEL ADARRAGA PINTO",  "1803", 52836,  1269000),
    @("73089307",   "GERMAN ANTONIO GARZON GOMEZ",     "1803", 42401,  1018372)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $ws.Cells.Item($row, 3).Value = $entry[0]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $entry[1]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $entry[2]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $entry[3]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value = $entry[4]   # G - Salario Basico
}
